# Apply GitHub-Actions crypto price refresh (Thu Sep 28 17:25:02 UTC 2023).
# Updates the Price (D) and Volume(1h) (E) columns across the existing rows,
# the row 41/42 swap (PaxDollar <-> MXToken), and the row 49-51 shuffle that
# inserts BabyDogeCoin (pushing Algorand/Cronos down a row and dropping the
# former EnergySwap row off the bottom of the list).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.105.87'
$ws.Range("E2").Value = '  +3.25%  '

# Row 3
$ws.Range("D3").Value = '1.659.05'

# Row 4
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").Value = '''215.41'
$ws.Range("E5").Value = '  +1.33%  '

# Row 6
$ws.Range("E6").Value = '  +1.15%  '

# Row 7
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  -0.14%  '

# Row 8
$ws.Range("E8").Value = '  +2.22%  '

# Row 9
$ws.Range("D9").Value = '''0.0614'
$ws.Range("E9").Value = '  +1.43%  '

# Row 10
$ws.Range("D10").Value = '''19.52'
$ws.Range("E10").Value = '  +3.05%  '

# Row 11
$ws.Range("E11").Value = '  +0.85%  '

# Row 12
$ws.Range("D12").Value = '1.893.84'
$ws.Range("E12").Value = '  +3.82%  '

# Row 13
$ws.Range("D13").Value = '1.660.61'
$ws.Range("E13").Value = '  +3.72%  '

# Row 14
$ws.Range("E14").Value = '  +2.03%  '

# Row 16
$ws.Range("D16").Value = '''64.93'
$ws.Range("E16").Value = '  +1.94%  '

# Row 17
$ws.Range("D17").Value = '''240.77'
$ws.Range("E17").Value = '  +5.86%  '

# Row 18
$ws.Range("D18").Value = '27.121.79'
$ws.Range("E18").Value = '  +3.31%  '

# Row 19
$ws.Range("D19").Value = '''7.86'
$ws.Range("E19").Value = '  +4.15%  '

# Row 20
$ws.Range("E20").Value = '  +1.39%  '

# Row 21
$ws.Range("E21").Value = '  -0.04%  '

# Row 22
$ws.Range("E22").Value = '  +5.11%  '

# Row 23
$ws.Range("D23").Value = '''2.27'
$ws.Range("E23").Value = '  +3.78%  '

# Row 24
$ws.Range("E24").Value = '  +3.96%  '

# Row 25
$ws.Range("D25").Value = '''145.93'
$ws.Range("E25").Value = '  +0.09%  '

# Row 26
$ws.Range("E26").Value = '  -0.06%  '

# Row 27
$ws.Range("E27").Value = '  +2.71%  '

# Row 28
$ws.Range("E28").Value = '  +1.18%  '

# Row 29
$ws.Range("D29").Value = '''15.86'
$ws.Range("E29").Value = '  +3.22%  '

# Row 30
$ws.Range("E30").Value = '  +1.06%  '

# Row 31
$ws.Range("E31").Value = '  +1.13%  '

# Row 32
$ws.Range("D32").Value = '1.526.05'
$ws.Range("E32").Value = '  +5.65%  '

# Row 33
$ws.Range("D33").Value = '''3.30'
$ws.Range("E33").Value = '  +3.12%  '

# Row 34
$ws.Range("E34").Value = '  +3.56%  '

# Row 35
$ws.Range("D35").Value = '''1.59'
$ws.Range("E35").Value = '  +8.10%  '

# Row 36
$ws.Range("E36").Value = '  -0.20%  '

# Row 37
$ws.Range("E37").Value = '  +1.41%  '

# Row 38
$ws.Range("D38").Value = '''0.892'
$ws.Range("E38").Value = '  +8.80%  '

# Row 39
$ws.Range("E39").Value = '  +2.69%  '

# Row 40
$ws.Range("D40").Value = '''5.96'
$ws.Range("E40").Value = '  +3.30%  '

# Row 41
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '''2.35'
$ws.Range("E41").Value = '  +8.03%  '

# Row 42
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").Value = '''1.00'
$ws.Range("E42").Value = '  -0.07%  '

# Row 43
$ws.Range("D43").Value = '''66.37'
$ws.Range("E43").Value = '  +9.74%  '

# Row 44
$ws.Range("D44").Value = '1.799.81'
$ws.Range("E44").Value = '  +3.58%  '

# Row 45
$ws.Range("D45").Value = '''0.772'
$ws.Range("E45").Value = '  +2.36%  '

# Row 46
$ws.Range("D46").Value = '''0.915'
$ws.Range("E46").Value = '  -1.60%  '

# Row 47
$ws.Range("D47").Value = '''90.53'
$ws.Range("E47").Value = '  +3.38%  '

# Row 48
$ws.Range("E48").Value = '  +3.48%  '

# Row 49
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0103'
$ws.Range("E49").Value = '  -1.41%  '

# Row 50
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '''0.0977'
$ws.Range("E50").Value = '  +3.11%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '''0.0502'
$ws.Range("E51").Value = '  +0.64%  '
